# Auto-generated: update cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.554.52'
$ws.Range('E2').Value = '  -1.94%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.627.53'
$ws.Range('E3').Value = '  -1.80%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.16'
$ws.Range('E5').Value = '  -3.58%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '156.64'
$ws.Range('E6').Value = '  -0.65%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').Value = '  +5.54%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('E9').Value = '  -5.25%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.80'
$ws.Range('E10').Value = '  -1.06%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.390'
$ws.Range('E11').Value = '  -2.61%  '

$ws.Range('E12').Value = '  +0.26%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '28.57'
$ws.Range('E13').Value = '  -2.08%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000187'
$ws.Range('E14').Value = '  -7.13%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.097.49'
$ws.Range('E15').Value = '  -1.84%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.410.64'
$ws.Range('E16').Value = '  -1.92%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.613.57'
$ws.Range('E17').Value = '  -2.06%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.28'
$ws.Range('E18').Value = '  -4.04%  '

$ws.Range('E19').Value = '  -2.75%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.38'
$ws.Range('E20').Value = '  -2.05%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '345.97'
$ws.Range('E21').Value = '  -1.91%  '

$ws.Range('E22').Value = '  -0.19%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.95'
$ws.Range('E23').Value = '  -2.39%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.0000112'
$ws.Range('E24').Value = '  -4.40%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.73'
$ws.Range('E25').Value = '  +3.48%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.39'
$ws.Range('E26').Value = '  -3.79%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.56'
$ws.Range('E27').Value = '  -2.70%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '556.63'
$ws.Range('E28').Value = '  +3.88%  '

$ws.Range('E29').Value = '  -2.35%  '

$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.00'
$ws.Range('E30').Value = '  -0.58%  '

$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('E32').Value = '  -3.25%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.74'
$ws.Range('E33').Value = '  -2.35%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.44'
$ws.Range('E34').Value = '  -1.13%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.36'
$ws.Range('E35').Value = '  -2.40%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.412'
$ws.Range('E36').Value = '  -2.64%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.05'
$ws.Range('E37').Value = '  -3.01%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.06%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.94'
$ws.Range('E39').Value = '  -0.40%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '151.61'
$ws.Range('E40').Value = '  -4.31%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.02%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.45'
$ws.Range('E42').Value = '  +2.96%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '158.45'
$ws.Range('E43').Value = '  -3.05%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.99'
$ws.Range('E44').Value = '  -2.95%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0602'
$ws.Range('E45').Value = '  -2.20%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.92'
$ws.Range('E46').Value = '  +0.62%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.635'
$ws.Range('E47').Value = '  -1.02%  '

$ws.Range('E48').Value = '  +2.49%  '

$ws.Range('E49').Value = '  -3.47%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '19.20'
$ws.Range('E50').Value = '  -4.47%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0240'
$ws.Range('E51').Value = '  -6.30%  '
